$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A2:F2").Copy($ws.Range("A3:F3"))
$ws.Range("A2:F2").Copy($ws.Range("A4:F4"))
$ws.Range("A2:F2").Copy($ws.Range("A5:F5"))
$ws.Range("A2:F2").Copy($ws.Range("A6:F6"))

$ws.Range("A5").Activate()
$ws.Range("A1:F6").Select()
